$d = $word.ActiveDocument

function Get-ParagraphIndexForRange($rng) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($rng.Start -ge $p.Range.Start -and $rng.Start -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Insert a new centered, bold "Nhóm 10" paragraph (Times New Roman, 16pt,
#    lang nl-NL) right before the three blank paragraphs that precede the
#    author's name ("VÕ ĐĂNG QUANG") on the cover page.
# ---------------------------------------------------------------------------

$quangFind = $d.Content
$quangFind.Find.ClearFormatting()
$quangFind.Find.Execute("VÕ ĐĂNG QUANG") | Out-Null
$quangIndex = Get-ParagraphIndexForRange $quangFind

$insertIndex = $quangIndex - 3
$targetPara = $d.Paragraphs.Item($insertIndex)
$targetPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($insertIndex)
$newRange = $newPara.Range

$groupXml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:before="80" w:after="80" w:line="360" w:lineRule="auto"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:b/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:lang w:val="nl-NL"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:b/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
                <w:lang w:val="nl-NL"/>
              </w:rPr>
              <w:t>Nh&#243;m 10</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newRange.InsertXML($groupXml)

# ---------------------------------------------------------------------------
# 2) Remove one of the two consecutive blank paragraphs that sit between
#    "BÁO CÁO CÁ NHÂN" and "CÁC CÔNG NGHỆ LẬP TRÌNH HIỆN ĐẠI" - specifically
#    the one immediately preceding the latter (carrying the left tab at 3930
#    and <w:b/>), keeping the earlier blank paragraph (carrying <w:bCs/>).
# ---------------------------------------------------------------------------

$techFind = $d.Content
$techFind.Find.ClearFormatting()
$techFind.Find.Execute("CÁC CÔNG NGHỆ LẬP TRÌNH HIỆN ĐẠI") | Out-Null
$techIndex = Get-ParagraphIndexForRange $techFind

$delIndex = $techIndex - 1
$delPara = $d.Paragraphs.Item($delIndex)
$delPara.Range.Delete()
